$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by duplicating "2021-Q4" (keeps the same
#    header/column layout + cell styling) and placing it right after
#    "2021-Q4" (i.e. right before "总计").
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# The template only has 3 data rows (rows 2-4); we need 7 (rows 2-8), so
# stretch the row-index column's style (s=2, bold/centered/bordered) down
# to the extra rows by copying the format of row 4's A cell.
$new.Range("A4").Copy()
$new.Range("A5:A8").PasteSpecial(-4122)

# Fund codes (column B) and the numeric-looking text fields (columns D-G)
# must stay text (e.g. leading zeros in fund codes, fixed decimal strings)
# instead of being auto-coerced to numbers, so force a Text format first.
$new.Range("B2:B8").NumberFormat = "@"
$new.Range("D2:G8").NumberFormat = "@"

$data = @(
  @("240004", "华宝动力组合混合", "10.49", "79.23", "3.34", "0.3504", 5),
  @("501029", "华宝标普中国A股红利机会指数（LOF）A", "13.19", "94.39", "1.98", "0.2612", 3),
  @("002906", "南方中证500量化增强股票A", "7.21", "92.26", "1.12", "0.0808", 5),
  @("005562", "创金合信中证红利低波动指数C", "1.85", "94.46", "2.06", "0.0381", 9),
  @("512890", "华泰柏瑞中证红利低波动ETF", "1.36", "99.24", "2.17", "0.0295", 9),
  @("005561", "创金合信中证红利低波动指数A", "1.22", "94.46", "2.06", "0.0251", 9),
  @("002907", "南方中证500量化增强股票C", "1.36", "92.26", "1.12", "0.0152", 5)
)

$r = 2
$idx = 0
foreach ($row in $data) {
  $new.Cells.Item($r, 1).Value = $idx
  $new.Cells.Item($r, 2).Value = $row[0]
  $new.Cells.Item($r, 3).Value = $row[1]
  $new.Cells.Item($r, 4).Value = $row[2]
  $new.Cells.Item($r, 5).Value = $row[3]
  $new.Cells.Item($r, 6).Value = $row[4]
  $new.Cells.Item($r, 7).Value = $row[5]
  $new.Cells.Item($r, 8).Value = $row[6]
  $r = $r + 1
  $idx = $idx + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new row right under the
#    header for the 2022-Q1 summary, pushing the older quarters down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.8

# The old rows (now shifted down one) keep their original index numbering in
# column A (0,1,2,3,4); renumber them 1,2,3,4,5 to stay sequential after the
# new row 0 above them.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(7, 1).Value = 5

Write-Output "2022-Q1 sheet + 总计 summary row added"
